# Delete the "~TFM_MIG" block (rows 9-11) from the INS sheet, which shifts
# everything below it up by three rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

$ws.Rows("9:11").Delete()

# Update the selection to reflect the new active cell (B9) as recorded
# in the saved workbook.
$ws.Range("B9").Select() | Out-Null
